# Fix a bug in image wall - greatly reduce CPU.
#
# Changes applied:
#  1. C24: append note about the AssistantHorizontalView button to the
#     existing "Update backend to store more images." comment.
#  2. B24: was a literal 1, now computed by the formula =1+3 (evaluates to 4).
#     B42 (=SUM(B2:B24)) recalculates automatically from 34 to 37.
#  3. Selection moves to B24 (was B42) to reflect where the edit happened.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the note text in C24 (shared string reused in place).
$ws.Range("C24").Value = "Update backend to store more images. Button on AssistantHorizontalView works now."

# 2. Replace the literal hours value in B24 with a formula; B42's SUM
#    picks up the new total automatically on recalculation.
$ws.Range("B24").Formula = "=1+3"

# 3. Move the selection/active cell to B24.
$ws.Range("B24").Select()
